$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 29 de Abril de 2020 a las 17:52"

# Row 4
$ws.Cells.Item(4, 2).Value = 1039501
$ws.Cells.Item(4, 3).Value = 3736
$ws.Cells.Item(4, 5).Value = 836925
$ws.Cells.Item(4, 6).Value = 19106
$ws.Cells.Item(4, 7).Value = 193
$ws.Cells.Item(4, 8).Value = 59459

# Row 9
$ws.Cells.Item(9, 5).Value = 33705
$ws.Cells.Item(9, 7).Value = 60
$ws.Cells.Item(9, 8).Value = 6374

# Row 15
$ws.Cells.Item(15, 2).Value = 50373
$ws.Cells.Item(15, 3).Value = 347
$ws.Cells.Item(15, 5).Value = 28279
$ws.Cells.Item(15, 7).Value = 45
$ws.Cells.Item(15, 8).Value = 2904

# Row 36
$ws.Cells.Item(36, 5).Value = 7724
$ws.Cells.Item(36, 7).Value = 22
$ws.Cells.Item(36, 8).Value = 685

# Row 44
$ws.Cells.Item(44, 6).Value = 40

# Row 57
$ws.Cells.Item(57, 1).Value = "Argelia"
$ws.Cells.Item(57, 2).Value = 3848
$ws.Cells.Item(57, 3).Value = 199
$ws.Cells.Item(57, 4).Value = 1702
$ws.Cells.Item(57, 5).Value = 1702
$ws.Cells.Item(57, 6).Value = 22
$ws.Cells.Item(57, 7).Value = 7
$ws.Cells.Item(57, 8).Value = 444

# Row 58
$ws.Cells.Item(58, 1).Value = "Moldavia"
$ws.Cells.Item(58, 2).Value = 3771
$ws.Cells.Item(58, 3).Value = 133
$ws.Cells.Item(58, 4).Value = 1114
$ws.Cells.Item(58, 5).Value = 2550
$ws.Cells.Item(58, 6).Value = 212
$ws.Cells.Item(58, 7).Value = 4
$ws.Cells.Item(58, 8).Value = 107

# Row 59
$ws.Cells.Item(59, 1).Value = "Luxemburgo"
$ws.Cells.Item(59, 2).Value = 3741
$ws.Cells.Item(59, 3).Value = 0
$ws.Cells.Item(59, 4).Value = 3123
$ws.Cells.Item(59, 5).Value = 529
$ws.Cells.Item(59, 6).Value = 19
$ws.Cells.Item(59, 7).Value = 0
$ws.Cells.Item(59, 8).Value = 89

# Row 60
$ws.Cells.Item(60, 1).Value = "Kuwait"
$ws.Cells.Item(60, 2).Value = 3740
$ws.Cells.Item(60, 3).Value = 300
$ws.Cells.Item(60, 4).Value = 1389
$ws.Cells.Item(60, 5).Value = 2327
$ws.Cells.Item(60, 6).Value = 66
$ws.Cells.Item(60, 7).Value = 1
$ws.Cells.Item(60, 8).Value = 24

# Row 65
$ws.Cells.Item(65, 2).Value = 2576
$ws.Cells.Item(65, 3).Value = 10
$ws.Cells.Item(65, 5).Value = 1860
$ws.Cells.Item(65, 6).Value = 41
$ws.Cells.Item(65, 7).Value = 1
$ws.Cells.Item(65, 8).Value = 139

# Row 80
$ws.Cells.Item(80, 1).Value = "Cuba"
$ws.Cells.Item(80, 2).Value = 1467
$ws.Cells.Item(80, 3).Value = 30
$ws.Cells.Item(80, 4).Value = 617
$ws.Cells.Item(80, 5).Value = 792
$ws.Cells.Item(80, 6).Value = 14
$ws.Cells.Item(80, 7).Value = 0
$ws.Cells.Item(80, 8).Value = 58

# Row 81
$ws.Cells.Item(81, 1).Value = "Bulgaria"
$ws.Cells.Item(81, 2).Value = 1447
$ws.Cells.Item(81, 3).Value = 48
$ws.Cells.Item(81, 4).Value = 243
$ws.Cells.Item(81, 5).Value = 1140
$ws.Cells.Item(81, 6).Value = 38
$ws.Cells.Item(81, 7).Value = 6
$ws.Cells.Item(81, 8).Value = 64

# Row 82
$ws.Cells.Item(82, 1).Value = "Republica de Macedonia"
$ws.Cells.Item(82, 2).Value = 1442
$ws.Cells.Item(82, 3).Value = 21
$ws.Cells.Item(82, 4).Value = 627
$ws.Cells.Item(82, 5).Value = 742
$ws.Cells.Item(82, 6).Value = 13
$ws.Cells.Item(82, 7).Value = 2
$ws.Cells.Item(82, 8).Value = 73

# Row 104
$ws.Cells.Item(104, 4).Value = 136
$ws.Cells.Item(104, 5).Value = 484

# Row 111
$ws.Cells.Item(111, 1).Value = "Mali"
$ws.Cells.Item(111, 2).Value = 482
$ws.Cells.Item(111, 3).Value = 58
$ws.Cells.Item(111, 4).Value = 129
$ws.Cells.Item(111, 5).Value = 328
$ws.Cells.Item(111, 6).Value = 0
$ws.Cells.Item(111, 7).Value = 1
$ws.Cells.Item(111, 8).Value = 25

# Row 112
$ws.Cells.Item(112, 1).Value = "Tanzania"
$ws.Cells.Item(112, 2).Value = 480
$ws.Cells.Item(112, 3).Value = 181
$ws.Cells.Item(112, 4).Value = 167
$ws.Cells.Item(112, 5).Value = 297
$ws.Cells.Item(112, 6).Value = 7
$ws.Cells.Item(112, 7).Value = 6
$ws.Cells.Item(112, 8).Value = 16

# Row 113
$ws.Cells.Item(113, 1).Value = "Malta"
$ws.Cells.Item(113, 2).Value = 463
$ws.Cells.Item(113, 3).Value = 5
$ws.Cells.Item(113, 4).Value = 339
$ws.Cells.Item(113, 5).Value = 120
$ws.Cells.Item(113, 6).Value = 1

# Row 114
$ws.Cells.Item(114, 1).Value = "Mayotte"
$ws.Cells.Item(114, 2).Value = 460
$ws.Cells.Item(114, 4).Value = 235
$ws.Cells.Item(114, 5).Value = 221
$ws.Cells.Item(114, 6).Value = 4
$ws.Cells.Item(114, 8).Value = 4

# Row 115
$ws.Cells.Item(115, 1).Value = "Jordania"
$ws.Cells.Item(115, 2).Value = 451
$ws.Cells.Item(115, 3).Value = 2
$ws.Cells.Item(115, 4).Value = 356
$ws.Cells.Item(115, 5).Value = 87
$ws.Cells.Item(115, 6).Value = 5
$ws.Cells.Item(115, 8).Value = 8

# Row 116
$ws.Cells.Item(116, 1).Value = "Taiwan"
$ws.Cells.Item(116, 2).Value = 429
$ws.Cells.Item(116, 4).Value = 311
$ws.Cells.Item(116, 5).Value = 112
$ws.Cells.Item(116, 8).Value = 6

# Row 121
$ws.Cells.Item(121, 1).Value = "Sudan"
$ws.Cells.Item(121, 2).Value = 375
$ws.Cells.Item(121, 3).Value = 57
$ws.Cells.Item(121, 4).Value = 32
$ws.Cells.Item(121, 5).Value = 315
$ws.Cells.Item(121, 7).Value = 3
$ws.Cells.Item(121, 8).Value = 28

# Row 122
$ws.Cells.Item(122, 1).Value = "Estado de Palestina"
$ws.Cells.Item(122, 2).Value = 344
$ws.Cells.Item(122, 3).Value = 1
$ws.Cells.Item(122, 4).Value = 71
$ws.Cells.Item(122, 5).Value = 271
$ws.Cells.Item(122, 6).Value = 0
$ws.Cells.Item(122, 8).Value = 2

# Row 123
$ws.Cells.Item(123, 1).Value = "Mauricio"
$ws.Cells.Item(123, 2).Value = 334
$ws.Cells.Item(123, 4).Value = 306
$ws.Cells.Item(123, 5).Value = 18

# Row 124
$ws.Cells.Item(124, 1).Value = "Venezuela"
$ws.Cells.Item(124, 2).Value = 329
$ws.Cells.Item(124, 4).Value = 142
$ws.Cells.Item(124, 5).Value = 177
$ws.Cells.Item(124, 6).Value = 3
$ws.Cells.Item(124, 8).Value = 10

# Row 125
$ws.Cells.Item(125, 1).Value = "Montenegro"
$ws.Cells.Item(125, 2).Value = 322
$ws.Cells.Item(125, 3).Value = 1
$ws.Cells.Item(125, 4).Value = 203
$ws.Cells.Item(125, 5).Value = 112
$ws.Cells.Item(125, 6).Value = 7
$ws.Cells.Item(125, 8).Value = 7

# Row 161
$ws.Cells.Item(161, 1).Value = "Islas Caimanes"
$ws.Cells.Item(161, 3).Value = 3
$ws.Cells.Item(161, 4).Value = 10
$ws.Cells.Item(161, 5).Value = 62
$ws.Cells.Item(161, 6).Value = 3

# Row 162
$ws.Cells.Item(162, 1).Value = "Guinea-Bisau"
$ws.Cells.Item(162, 2).Value = 73
$ws.Cells.Item(162, 4).Value = 18
$ws.Cells.Item(162, 5).Value = 54

# Row 163
$ws.Cells.Item(163, 1).Value = "Suazilandia"
$ws.Cells.Item(163, 2).Value = 71
$ws.Cells.Item(163, 5).Value = 60
$ws.Cells.Item(163, 6).Value = 0

# Row 183
$ws.Cells.Item(183, 4).Value = 5
$ws.Cells.Item(183, 5).Value = 17
